# Applies the data refresh captured in the commit "Update gh-pages to
# output generated at 456a3b4" to the workbook.
#
# Sheet "展览" (exhibitions): a brand-new entry ("北京·井上直久の依巴拉度世界")
# is published into row 2, the two previously-top entries shift down one
# row (row2->row3, row3->row4), the old row4 entry ("北京·第五元素二次元音乐
# 天堂（取消）") drops off the list, and many of the remaining rows get an
# updated "想去人数" (interested-count) and occasionally a refreshed cover
# image. Column A (index) and column B (date) are left untouched.
#
# Sheet "演出" (performances): two small 想去人数 count bumps.
#
# Sheet "全部类型" (all types): the same 想去人数 / cover-image refresh as
# sheet "展览", but this sheet's rows were not shifted.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览 (exhibitions)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

# New row 2: 北京·井上直久の依巴拉度世界
$ws1.Range("C2").Value = "北京·井上直久の依巴拉度世界"
$ws1.Range("D2").Value = "798艺术区东街 美仑美术馆"
$ws1.Range("E2").Value = "2024.06.15 10:00-08.08 19:00"
$ws1.Range("F2").Value = 0
$ws1.Range("G2").Value = 19.9
$ws1.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=87162"
$ws1.Range("I2").Value = "//i2.hdslb.com/bfs/openplatform/202406/JRwxed341718092866573.jpeg"

# Row 3 becomes what row 2 used to be: 北京·代号鸢ONLY·女仆主题日
$ws1.Range("C3").Value = "北京·代号鸢ONLY·女仆主题日"
$ws1.Range("D3").Value = "垡头街道双合北街10号楼(双合地铁站C口旁) 合憬荟"
$ws1.Range("E3").Value = "2024.06.15 09:00-06.16 17:00"
$ws1.Range("F3").Value = 238
$ws1.Range("G3").Value = 88
$ws1.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=86352"
$ws1.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202405/sHSYcfWj1716704297772.jpeg"

# Row 4 becomes what row 3 used to be: 北京·恋与深空only（取消）
# (G4 keeps its existing "不可售" value, so it is left alone.)
$ws1.Range("C4").Value = "北京·恋与深空only（取消）"
$ws1.Range("D4").Value = "太平庄中街西端 北京天通苑黄河京都会议中心"
$ws1.Range("E4").Value = "2024.06.15 10:00-06.15 17:00"
$ws1.Range("F4").Value = 375
$ws1.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=84729"
$ws1.Range("I4").Value = "//i2.hdslb.com/bfs/openplatform/202404/VWBsXunu1713865945134.jpeg"

# Old row 4 (北京·第五元素二次元音乐天堂（取消）) drops out of the list;
# rows 5 and below keep their position, only their 想去人数 (column F)
# -- and for row 38 also the cover image -- are refreshed below.

$sheet1Updates = @{
    5  = 25
    6  = 93
    8  = 384
    9  = 4648
    10 = 4648
    11 = 127
    14 = 605
    15 = 4149
    16 = 160
    17 = 161
    18 = 48
    19 = 203
    20 = 3424
    24 = 2988
    25 = 124
    26 = 124
    27 = 7
    28 = 143
    29 = 180
    31 = 72
    32 = 49
    36 = 5380
    37 = 752
    38 = 384
    41 = 31
    42 = 1072
    43 = 447
    45 = 1942
    47 = 53
    48 = 691
    49 = 828
}

foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Refreshed cover image for row 38 (万游引力国潮动漫嘉年华s8)
$ws1.Range("I38").Value = "//i1.hdslb.com/bfs/openplatform/202406/0kySwWBG1718096478563.jpeg"

# ---------------------------------------------------------------------
# Sheet 2: 演出 (performances)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("演出")

$ws2.Cells.Item(15, 6).Value = 119
$ws2.Cells.Item(22, 6).Value = 722

# ---------------------------------------------------------------------
# Sheet 4: 全部类型 (all types)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$sheet4Updates = @{
    5  = 238
    6  = 25
    8  = 93
    10 = 384
    11 = 4648
    12 = 4648
    14 = 127
    19 = 605
    20 = 4149
    21 = 160
    22 = 161
    23 = 203
    24 = 3424
    25 = 2988
    26 = 124
    27 = 124
    28 = 143
    29 = 180
    35 = 119
    37 = 5380
    39 = 752
    40 = 384
    44 = 1072
    45 = 447
    47 = 1942
    48 = 53
    49 = 691
    50 = 828
}

foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}

# Refreshed cover image for row 40 (万游引力国潮动漫嘉年华s8)
$ws4.Range("I40").Value = "//i1.hdslb.com/bfs/openplatform/202406/0kySwWBG1718096478563.jpeg"
